$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '63.567.40'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -0.44%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.256.68'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +3.89%  '
$ws.Range("E4").Value = '  +0.01%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '596.56'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.25%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '141.78'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +1.83%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '3.251.15'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +3.96%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.518'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -1.37%  '
$ws.Range("E10").Value = '  +0.11%  '
$ws.Range("E11").Value = '  +1.45%  '
$ws.Range("E12").Value = '  +0.98%  '
$ws.Range("E13").Value = '  -1.65%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '34.46'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +0.14%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '3.790.25'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +3.86%  '
$ws.Range("E16").Value = '  -0.16%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '3.252.67'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +3.84%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '63.529.51'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -0.48%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '6.80'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +0.63%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '479.18'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -0.46%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '14.27'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -1.42%  '
$ws.Range("E22").Value = '  +3.91%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '7.98'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +3.65%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '83.84'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -4.37%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '13.33'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +2.36%  '
$ws.Range("E27").Value = '  +0.35%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '7.22'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +4.17%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '8.09'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -0.43%  '
$ws.Range("E30").Value = '  +4.56%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '27.73'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +2.24%  '
$ws.Range("E32").Value = '  -0.02%  '
$ws.Range("E33").Value = '  -3.32%  '
$ws.Range("E34").Value = '  -0.91%  '
$ws.Range("E35").Value = '  -0.75%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '5.94'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -0.99%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '52.79'
$ws.Range("D37").Style = "Normal"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.0₃0720'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -1.98%  '
$ws.Range("E39").Value = '  +0.36%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '423.46'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -0.86%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '2.998.10'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +4.42%  '
$ws.Range("E42").Value = '  -2.63%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '8.40'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +1.19%  '
$ws.Range("E44").Value = '  -7.71%  '
$ws.Range("E45").Value = '  +2.55%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.20'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +2.08%  '
$ws.Range("E47").Value = '  +0.09%  '
$ws.Range("E48").Value = '  -1.38%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '26.00'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +1.90%  '
$ws.Range("E50").Value = '  +0.47%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '122.78'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +1.89%  '
